$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (shifts existing rows down)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).RowHeight = 17

$ws.Cells.Item(4, 1).Value = "Luis D. Verde Arregoitia"
$ws.Cells.Item(4, 2).Value = "Large Language Model tools for R"
$ws.Cells.Item(4, 3).Value = "https://luisdva.github.io/llmsr-book/"

$ws.Range("A2").Select()
